$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B2" = 0.2245614035087719
    "C2" = 0.5087719298245614
    "J2" = 0.02807017543859649
    "P2" = 0.1473684210526316
    "S2" = 0.0912280701754386
    "B3" = 0.006535947712418301
    "C3" = 0.0392156862745098
    "J3" = 0.0130718954248366
    "P3" = 0.7450980392156863
    "S3" = 0.196078431372549
    "J4" = 0.1071428571428571
    "P4" = 0.7142857142857143
    "S4" = 0.1785714285714286
    "B6" = 0.07894736842105263
    "D6" = 0.02105263157894737
    "F6" = 0.06842105263157895
    "J6" = 0.1736842105263158
    "O6" = 0.04736842105263158
    "Q6" = 0.1052631578947368
    "R6" = 0.08947368421052632
    "S6" = 0.4157894736842105
    "B7" = 0.1151832460732984
    "D7" = 0.02094240837696335
    "E7" = 0.005235602094240838
    "F7" = 0.02617801047120419
    "J7" = 0.1413612565445026
    "O7" = 0.01047120418848168
    "Q7" = 0.193717277486911
    "R7" = 0.06806282722513089
    "S7" = 0.418848167539267
    "B8" = 0.1053921568627451
    "D8" = 0.01225490196078431
    "E8" = 0.002450980392156863
    "F8" = 0.05147058823529412
    "J8" = 0.1102941176470588
    "O8" = 0.02450980392156863
    "Q8" = 0.1691176470588235
    "R8" = 0.08333333333333333
    "S8" = 0.4411764705882353
    "B9" = 0.07058823529411765
    "F9" = 0.1058823529411765
    "J9" = 0.1529411764705882
    "O9" = 0.01176470588235294
    "Q9" = 0.1588235294117647
    "R9" = 0.07647058823529412
    "S9" = 0.4235294117647059
    "B10" = 0.1088850174216028
    "D10" = 0.01393728222996516
    "E10" = 0.0008710801393728223
    "F10" = 0.07926829268292683
    "J10" = 0.1332752613240418
    "O10" = 0.02264808362369338
    "Q10" = 0.1898954703832753
    "R10" = 0.08797909407665505
    "S10" = 0.3632404181184669
    "G11" = 0.1432926829268293
    "J11" = 0.1128048780487805
    "K11" = 0.2195121951219512
    "L11" = 0.4908536585365854
    "S11" = 0.03353658536585366
    "G12" = 0.7168674698795181
    "J12" = 0.2289156626506024
    "K12" = 0.01807228915662651
    "L12" = 0.01204819277108434
    "S12" = 0.02409638554216868
    "G13" = 0.7
    "J13" = 0.24
    "S13" = 0.06
    "F15" = 0.0187793427230047
    "H15" = 0.2065727699530517
    "I15" = 0.06572769953051644
    "J15" = 0.2816901408450704
    "K15" = 0.05633802816901409
    "M15" = 0.004694835680751174
    "O15" = 0.06103286384976526
    "S15" = 0.3051643192488263
    "F16" = 0.005747126436781609
    "H16" = 0.1781609195402299
    "I16" = 0.09195402298850575
    "J16" = 0.3563218390804598
    "K16" = 0.1494252873563219
    "M16" = 0.02873563218390805
    "O16" = 0.04597701149425287
    "S16" = 0.1436781609195402
    "F17" = 0.008086253369272238
    "H17" = 0.1455525606469003
    "I17" = 0.1078167115902965
    "J17" = 0.3908355795148248
    "K17" = 0.1293800539083558
    "M17" = 0.02425876010781671
    "O17" = 0.06199460916442048
    "S17" = 0.1320754716981132
    "F18" = 0.005649717514124294
    "H18" = 0.2542372881355932
    "I18" = 0.05649717514124294
    "J18" = 0.4067796610169492
    "K18" = 0.07909604519774012
    "M18" = 0.01129943502824859
    "O18" = 0.05084745762711865
    "S18" = 0.1355932203389831
    "F19" = 0.009159034138218152
    "H19" = 0.2014987510407993
    "I19" = 0.07660283097418817
    "J19" = 0.3621981681931724
    "K19" = 0.1248959200666112
    "M19" = 0.02830974188176519
    "N19" = 0.0008326394671107411
    "O19" = 0.07327227310574522
    "S19" = 0.1232306411323897
}

foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}

Write-Host "Applied $($changes.Count) cell updates"